# Populate the Mongolian localization QA sheet (open saved query file).
# Header row (A1:D1) + 9 rows of English/Mongolian translation pairs,
# duplicated twice (rows 2-10 and 11-19); "modified translation" / "reason
# of correction" columns (C, D) are left blank for every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'original English sentence'
$ws.Range('B1').Value = 'original translation'
$ws.Range('C1').Value = 'modified translation'
$ws.Range('D1').Value = 'reason of correction'

$ws.Range('A2').Value = 'Help fashion-challenged Danny pick a swell outfit for his big date with Jamie.'
$ws.Range('B2').Value = 'Загвар өмсөгч Дэннид Жэймитэй хийх том болзоондоо гоё хувцас сонгоход нь туслаарай.'
$ws.Range('C2').Value = ''
$ws.Range('D2').Value = ''

$ws.Range('A3').Value = 'Hit ''Stop'' to select the drawer containing the costume you want Danny to wear.'
$ws.Range('B3').Value = '"Зогс" дээр дарж Даннигийн өмсөхийг хүссэн хувцасны шүүгээг сонгоно уу.'
$ws.Range('C3').Value = ''
$ws.Range('D3').Value = ''

$ws.Range('A4').Value = 'Please don''t show me this dialogue again.'
$ws.Range('B4').Value = 'Энэ харилцан яриаг надад дахиж битгий үзүүлээрэй.'
$ws.Range('C4').Value = ''
$ws.Range('D4').Value = ''

$ws.Range('A5').Value = 'Avoid'
$ws.Range('B5').Value = 'Зайлсхий'
$ws.Range('C5').Value = ''
$ws.Range('D5').Value = ''

$ws.Range('A6').Value = 'Where''s Danny?'
$ws.Range('B6').Value = 'Данни хаана байна'
$ws.Range('C6').Value = ''
$ws.Range('D6').Value = ''

$ws.Range('A7').Value = 'Find'
$ws.Range('B7').Value = 'Хай'
$ws.Range('C7').Value = ''
$ws.Range('D7').Value = ''

$ws.Range('A8').Value = 'That miserable little punk Danny is trying to hide from me again. I''ll make it worth your while if you help me find him. But I''d rather avoid letting Jamie know I''m looking for Danny so let''s avoid her.'
$ws.Range('B8').Value = 'Тэр хөөрхийлөлтэй бяцхан панк Дэнни дахиад л надаас нуугдах гэж байна. Хэрэв та түүнийг олоход надад тусалбал би үүнийг үнэ цэнэтэй болгоно. Гэхдээ би Дэнниг хайж байгаагаа Жэймид мэдэгдэхээс зайлсхийсэн нь дээр, тиймээс түүнээс зайлсхийцгээе.'
$ws.Range('C8').Value = ''
$ws.Range('D8').Value = ''

$ws.Range('A9').Value = 'Click on the spots you think Danny might be hiding. Likely spots earn you coins.'
$ws.Range('B9').Value = 'Дэннигийн нуугдаж байгаа гэж бодож буй газрууд дээр дарна уу. Магадгүй цэгүүд танд зоос олох болно.'
$ws.Range('C9').Value = ''
$ws.Range('D9').Value = ''

$ws.Range('A10').Value = 'The round ends if you find Jamie'
$ws.Range('B10').Value = 'Хэрэв та Жэймиг олвол тойрог дуусна'
$ws.Range('C10').Value = ''
$ws.Range('D10').Value = ''

$ws.Range('A11').Value = 'Help fashion-challenged Danny pick a swell outfit for his big date with Jamie.'
$ws.Range('B11').Value = 'Загвар өмсөгч Дэннид Жэймитэй хийх том болзоондоо гоё хувцас сонгоход нь туслаарай.'
$ws.Range('C11').Value = ''
$ws.Range('D11').Value = ''

$ws.Range('A12').Value = 'Hit ''Stop'' to select the drawer containing the costume you want Danny to wear.'
$ws.Range('B12').Value = '"Зогс" дээр дарж Даннигийн өмсөхийг хүссэн хувцасны шүүгээг сонгоно уу.'
$ws.Range('C12').Value = ''
$ws.Range('D12').Value = ''

$ws.Range('A13').Value = 'Please don''t show me this dialogue again.'
$ws.Range('B13').Value = 'Энэ харилцан яриаг надад дахиж битгий үзүүлээрэй'
$ws.Range('C13').Value = ''
$ws.Range('D13').Value = ''

$ws.Range('A14').Value = 'Avoid'
$ws.Range('B14').Value = 'Зайлсхий'
$ws.Range('C14').Value = ''
$ws.Range('D14').Value = ''

$ws.Range('A15').Value = 'Where''s Danny?'
$ws.Range('B15').Value = 'Данни хаана байна'
$ws.Range('C15').Value = ''
$ws.Range('D15').Value = ''

$ws.Range('A16').Value = 'Find'
$ws.Range('B16').Value = 'Хай'
$ws.Range('C16').Value = ''
$ws.Range('D16').Value = ''

$ws.Range('A17').Value = 'That miserable little punk Danny is trying to hide from me again. I''ll make it worth your while if you help me find him. But I''d rather avoid letting Jamie know I''m looking for Danny so let''s avoid her.'
$ws.Range('B17').Value = 'Тэр хөөрхийлөлтэй бяцхан панк Дэнни дахиад л надаас нуугдах гэж байна. Хэрэв та түүнийг олоход надад тусалбал би үүнийг үнэ цэнэтэй болгоно. Гэхдээ би Дэнниг хайж байгаагаа Жэймид мэдэгдэхээс зайлсхийсэн нь дээр, тиймээс түүнээс зайлсхийцгээе.'
$ws.Range('C17').Value = ''
$ws.Range('D17').Value = ''

$ws.Range('A18').Value = 'Click on the spots you think Danny might be hiding. Likely spots earn you coins.'
$ws.Range('B18').Value = 'Дэннигийн нуугдаж байгаа гэж бодож буй газрууд дээр дарна уу. Магадгүй цэгүүд танд зоос олох болно.'
$ws.Range('C18').Value = ''
$ws.Range('D18').Value = ''

$ws.Range('A19').Value = 'The round ends if you find Jamie'
$ws.Range('B19').Value = 'Хэрэв та Жэймиг олвол тойрог дуусна'
$ws.Range('C19').Value = ''
$ws.Range('D19').Value = ''
